# Insert two new weekly price rows for "Haba" (Mercado Mayorista Lo Valledor
# de Santiago) ahead of the existing row 226, shifting the rest of the table
# down by two rows (226-300 -> 228-302), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 226 (pushes old rows 226..300 to 228..302)
$ws.Rows("226:227").Insert()

# --- New row 226 ---
$ws.Range("A226").Value = 6
$ws.Range("B226").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C226").Value = "Metropolitana"
$ws.Range("D226").Value = 44809
$ws.Range("E226").Value = 13
$ws.Range("F226").Value = 100112026
$ws.Range("G226").Value = "Haba"
$ws.Range("H226").Value = "Sin especificar"
$ws.Range("I226").Value = "Primera"
$ws.Range("J226").Value = 400
$ws.Range("K226").Value = 9000
$ws.Range("L226").Value = 10000
$ws.Range("M226").Value = 9425
$ws.Range("N226").Value = "`$/saco 25 kilos"
$ws.Range("O226").Value = "Provincia de Huasco"
$ws.Range("P226").Value = 377
$ws.Range("Q226").Value = 25
$ws.Range("R226").Value = "Hortaliza"

# --- New row 227 ---
$ws.Range("A227").Value = 6
$ws.Range("B227").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C227").Value = "Metropolitana"
$ws.Range("D227").Value = 44809
$ws.Range("E227").Value = 13
$ws.Range("F227").Value = 100112026
$ws.Range("G227").Value = "Haba"
$ws.Range("H227").Value = "Sin especificar"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 180
$ws.Range("K227").Value = 9000
$ws.Range("L227").Value = 10000
$ws.Range("M227").Value = 9444
$ws.Range("N227").Value = "`$/saco 25 kilos"
$ws.Range("O227").Value = "Provincia de Limarí"
$ws.Range("P227").Value = 378
$ws.Range("Q227").Value = 25
$ws.Range("R227").Value = "Hortaliza"
